$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Ingredients" title text from the merged A1:B1 header band
$ws.Range("A1").Value = $null

# Delete the blank spacer column A (IngrID numbers were in B, names in C).
# This shifts B -> A and C -> B.
$ws.Columns("A").Delete()

# Fix up the column headers on row 2: A2 = "IngrID", B2 = "IngrName"
$ws.Range("A2").Value = "IngrID"
$ws.Range("B2").Value = "IngrName"
